# gif updates, map bottom scroll update
#
# The "Column Name"/"Description" reference sheet gets two of its
# descriptions reworded (the underlying shared-string table shuffles as a
# side effect, but the only user-visible change is the text of B15 and
# B17), and the sheet's scroll/selection state is updated so the map
# bottom of the table (around row 18) is back in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reword the two descriptions (existing_solar_mw / suitable_rooftop_area_sqft)
$ws.Range("B15").Value = "Existing in front of the meter solar generation (MW)"
$ws.Range("B17").Value = "Total rooftop area (sqft) of institutional site considered suitable for solar PV"

# Update the view: scroll so row 10 is the top visible row, and move the
# active selection down to B18.
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2
